$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.071020364761353
$ws.Range("B1").Value = 3.357789516448975
$ws.Range("C1").Value = 2.835106611251831
$ws.Range("D1").Value = 2.12580132484436
$ws.Range("E1").Value = 1.246389150619507
